$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4084.4285
$ws.Range("J112").Value = 4084.4285
$ws.Range("L112").Value = 12253.2855
$ws.Range("N112").Value = -14469.2855
$ws.Range("H132").Value = 975.5484
$ws.Range("I132").Value = 960.069
$ws.Range("K132").Value = 2880.207
$ws.Range("M132").Value = -350.2069999999999
$ws.Range("H138").Value = 3486.818
$ws.Range("I138").Value = 6261.375
$ws.Range("J138").Value = 2870.25
$ws.Range("K138").Value = 18784.125
$ws.Range("L138").Value = 8610.75
$ws.Range("M138").Value = -13644.125
$ws.Range("N138").Value = -18890.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3221.4443
$ws.Range("I32").Value = 2628.4355
$ws.Range("J32").Value = 6898.1
$ws.Range("K32").Value = 2628.4355
$ws.Range("L32").Value = 6898.1
$ws.Range("M32").Value = -2341.4355
$ws.Range("N32").Value = -7472.1
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 1773.1428
$ws.Range("I122").Value = 1755.6923
$ws.Range("K122").Value = 5267.0769
$ws.Range("M122").Value = -2817.0769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2453.7
$ws.Range("I20").Value = 2417.8333
$ws.Range("J20").Value = 2507.5
$ws.Range("K20").Value = 2417.8333
$ws.Range("L20").Value = 2507.5
$ws.Range("M20").Value = -2170.8333
$ws.Range("N20").Value = -3001.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 26000
$ws.Range("J70").Value = 29000
$ws.Range("L70").Value = 29000
$ws.Range("N70").Value = -29630
$ws.Range("H73").Value = 26000
$ws.Range("J73").Value = 29000
$ws.Range("L73").Value = 29000
$ws.Range("N73").Value = -31184
$ws.Range("H94").Value = 1058.6
$ws.Range("I94").Value = 917.2
$ws.Range("J94").Value = 1200
$ws.Range("K94").Value = 917.2
$ws.Range("L94").Value = 1200
$ws.Range("M94").Value = -466.2
$ws.Range("N94").Value = -2102
$ws.Range("H107").Value = 551.2857
$ws.Range("I107").Value = 403.83334
$ws.Range("K107").Value = 403.83334
$ws.Range("M107").Value = 1516.16666
$ws.Range("H132").Value = 2733.879
$ws.Range("I132").Value = 2216.5789
$ws.Range("J132").Value = 3435.9285
$ws.Range("K132").Value = 6649.736699999999
$ws.Range("L132").Value = 10307.7855
$ws.Range("M132").Value = -4119.736699999999
$ws.Range("N132").Value = -15367.7855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 359.73914
$ws.Range("J5").Value = 337.35294
$ws.Range("L5").Value = 1012.05882
$ws.Range("N5").Value = -1236.05882
$ws.Range("H68").Value = 2235
$ws.Range("I68").Value = 1192.2142
$ws.Range("J68").Value = 2818.96
$ws.Range("K68").Value = 3576.6426
$ws.Range("L68").Value = 8456.880000000001
$ws.Range("M68").Value = -2765.6426
$ws.Range("N68").Value = -10078.88
$ws.Range("H71").Value = 2235
$ws.Range("I71").Value = 1192.2142
$ws.Range("J71").Value = 2818.96
$ws.Range("K71").Value = 10729.9278
$ws.Range("L71").Value = 25370.64
$ws.Range("M71").Value = -6673.927799999999
$ws.Range("N71").Value = -33482.64
$ws.Range("H107").Value = 1411.4423
$ws.Range("I107").Value = 1185.7142
$ws.Range("J107").Value = 1446.5555
$ws.Range("K107").Value = 3557.1426
$ws.Range("L107").Value = 4339.666499999999
$ws.Range("M107").Value = -1637.1426
$ws.Range("N107").Value = -8179.666499999999
$ws.Range("H120").Value = 166668180
$ws.Range("I120").Value = 166668180
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 500004540
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -499999702
$ws.Range("N120").ClearContents()
$ws.Range("H131").Value = 13909496
$ws.Range("J131").Value = 21789.47
$ws.Range("L131").Value = 65368.41
$ws.Range("N131").Value = -75448.41
$ws.Range("H135").Value = 359.73914
$ws.Range("J135").Value = 337.35294
$ws.Range("L135").Value = 3036.17646
$ws.Range("N135").Value = -8106.17646
$ws.Range("H137").Value = 3481.8462
$ws.Range("J137").Value = 5458.8
$ws.Range("L137").Value = 16376.4
$ws.Range("N137").Value = -26576.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2242
$ws.Range("I80").Value = 2241
$ws.Range("J80").Value = 2242.5
$ws.Range("K80").Value = 2241
$ws.Range("L80").Value = 2242.5
$ws.Range("M80").Value = -1243
$ws.Range("N80").Value = -4238.5
$ws.Range("H83").Value = 2242
$ws.Range("I83").Value = 2241
$ws.Range("J83").Value = 2242.5
$ws.Range("K83").Value = 11205
$ws.Range("L83").Value = 11212.5
$ws.Range("M83").Value = -6213
$ws.Range("N83").Value = -21196.5
$ws.Range("H97").Value = 1378.4615
$ws.Range("I97").Value = 1035.8889
$ws.Range("K97").Value = 1035.8889
$ws.Range("M97").Value = -539.8888999999999
$ws.Range("H102").Value = 3808.4443
$ws.Range("I102").Value = 3909.2856
$ws.Range("K102").Value = 3909.2856
$ws.Range("M102").Value = -2287.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2747.125
$ws.Range("I68").Value = 2425.2856
$ws.Range("K68").Value = 2425.2856
$ws.Range("M68").Value = -1676.2856
$ws.Range("H71").Value = 2747.125
$ws.Range("I71").Value = 2425.2856
$ws.Range("K71").Value = 12126.428
$ws.Range("M71").Value = -8382.428
$ws.Range("H93").Value = 759.5333000000001
$ws.Range("I93").Value = 778.0714
$ws.Range("K93").Value = 778.0714
$ws.Range("M93").Value = 469.9286
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680
$ws.Range("H132").Value = 3708.4348
$ws.Range("I132").Value = 1024.6666
$ws.Range("J132").Value = 6636.1816
$ws.Range("K132").Value = 3073.9998
$ws.Range("L132").Value = 19908.5448
$ws.Range("M132").Value = -543.9998000000001
$ws.Range("N132").Value = -24968.5448
$ws.Range("H136").Value = 3778.6667
$ws.Range("I136").Value = 2611.6365
$ws.Range("K136").Value = 7834.9095
$ws.Range("M136").Value = -5284.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 10000
$ws.Range("J112").Value = 10000
$ws.Range("L112").Value = 10000
$ws.Range("N112").Value = -12954
$ws.Range("H132").Value = 571.1818
$ws.Range("I132").Value = 571.1818
$ws.Range("K132").Value = 1713.5454
$ws.Range("M132").Value = 816.4546
